$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").ClearContents() | Out-Null
$ws.Range("N19").ClearContents() | Out-Null

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").ClearContents() | Out-Null
$ws.Range("N46").ClearContents() | Out-Null

$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").ClearContents() | Out-Null
$ws.Range("N60").ClearContents() | Out-Null

$ws.Range("H92").Value = 2098.7917
$ws.Range("I92").Value = 202.4375
$ws.Range("J92").Value = 5891.5
$ws.Range("K92").Value = 202.4375
$ws.Range("L92").Value = 5891.5
$ws.Range("M92").Value = 1045.5625
$ws.Range("N92").Value = -8387.5

$ws.Range("H96").Value = 6494572
$ws.Range("I96").Value = 17857476
$ws.Range("J96").Value = 1484.1428
$ws.Range("K96").Value = 53572428
$ws.Range("L96").Value = 4452.428400000001
$ws.Range("M96").Value = -53571055
$ws.Range("N96").Value = -7198.428400000001

$ws.Range("H99").Value = 509.92856
$ws.Range("I99").Value = 509.92856
$ws.Range("K99").Value = 1529.78568
$ws.Range("M99").Value = -31.78567999999996

$ws.Range("H101").Value = 280
$ws.Range("I101").Value = 292.22223
$ws.Range("K101").Value = 876.66669
$ws.Range("M101").Value = 745.33331

$ws.Range("H106").Value = 3568.75
$ws.Range("I106").Value = 3568.75
$ws.Range("K106").Value = 3568.75
$ws.Range("M106").Value = -2937.75

$ws.Range("H113").Value = 5226.727
$ws.Range("J113").Value = 6415.8335
$ws.Range("L113").Value = 6415.8335
$ws.Range("N113").Value = -12923.8335

$ws.Range("H125").Value = 75766.5
$ws.Range("J125").Value = 1017.5
$ws.Range("L125").Value = 9157.5
$ws.Range("N125").Value = -14077.5

$ws.Range("H127").Value = 2860.5334
$ws.Range("J127").Value = 4881.25
$ws.Range("L127").Value = 14643.75
$ws.Range("N127").Value = -24563.75

$ws.Range("H132").Value = 34010.453
$ws.Range("I132").Value = 37071.8
$ws.Range("K132").Value = 111215.4
$ws.Range("M132").Value = -108685.4

$ws.Range("H138").Value = 27424.15
$ws.Range("I138").Value = 1727.6818
$ws.Range("K138").Value = 5183.0454
$ws.Range("M138").Value = -43.04539999999997

$ws.Range("H141").Value = 1372.2727
$ws.Range("I141").Value = 1261.5
$ws.Range("K141").Value = 3784.5
$ws.Range("M141").Value = 1395.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 3225.1667
$ws.Range("I39").Value = 3225.1667
$ws.Range("K39").Value = 3225.1667
$ws.Range("M39").Value = -2705.1667

$ws.Range("H97").Value = 1629.8636
$ws.Range("I97").Value = 1242.091
$ws.Range("J97").Value = 2017.6364
$ws.Range("K97").Value = 1242.091
$ws.Range("L97").Value = 2017.6364
$ws.Range("M97").Value = -746.0909999999999
$ws.Range("N97").Value = -3009.6364

$ws.Range("H122").Value = 1378.7778
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 5000
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents() | Out-Null

$ws.Range("H107").Value = 6077.5356
$ws.Range("I107").Value = 6087.12
$ws.Range("J107").Value = 5997.6665
$ws.Range("K107").Value = 6087.12
$ws.Range("L107").Value = 5997.6665
$ws.Range("M107").Value = -4167.12
$ws.Range("N107").Value = -9837.666499999999

$ws.Range("H134").Value = 2949.318
$ws.Range("I134").Value = 2254.7334
$ws.Range("K134").Value = 6764.2002
$ws.Range("M134").Value = -4229.2002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 180

$ws.Range("H31").Value = 3705367.2
$ws.Range("I31").Value = 5000961
$ws.Range("J31").Value = 3671.4285
$ws.Range("K31").Value = 5000961
$ws.Range("L31").Value = 3671.4285
$ws.Range("M31").Value = -5000666
$ws.Range("N31").Value = -4261.4285

$ws.Range("H34").Value = 3705367.2
$ws.Range("I34").Value = 5000961
$ws.Range("J34").Value = 3671.4285
$ws.Range("K34").Value = 5000961
$ws.Range("L34").Value = 3671.4285
$ws.Range("M34").Value = -5000759
$ws.Range("N34").Value = -4075.4285

$ws.Range("H58").Value = 14121.115
$ws.Range("I58").Value = 1369.85
$ws.Range("K58").Value = 1369.85
$ws.Range("M58").Value = -1166.85

$ws.Range("H62").Value = 6908.6313
$ws.Range("I62").Value = 6742.1875
$ws.Range("K62").Value = 6742.1875
$ws.Range("M62").Value = -6118.1875

$ws.Range("H65").Value = 6908.6313
$ws.Range("I65").Value = 6742.1875
$ws.Range("K65").Value = 33710.9375
$ws.Range("M65").Value = -30590.9375

$ws.Range("H99").Value = 6884.5713
$ws.Range("I99").Value = 5623.125
$ws.Range("K99").Value = 5623.125
$ws.Range("M99").Value = -4125.125

$ws.Range("H107").Value = 724.12
$ws.Range("I107").Value = 666.5625
$ws.Range("K107").Value = 666.5625
$ws.Range("M107").Value = 1253.4375

$ws.Range("H122").Value = 1581.7693
$ws.Range("I122").Value = 1627.1052
$ws.Range("K122").Value = 4881.3156
$ws.Range("M122").Value = -2431.3156

$ws.Range("H126").Value = 6884.5713
$ws.Range("I126").Value = 5623.125
$ws.Range("K126").Value = 16869.375
$ws.Range("M126").Value = -14399.375

$ws.Range("H132").Value = 49077.285
$ws.Range("I132").Value = 63344.625
$ws.Range("K132").Value = 190033.875
$ws.Range("M132").Value = -187503.875

$ws.Range("H134").Value = 2146.3157
$ws.Range("I134").Value = 1810.5883
$ws.Range("K134").Value = 5431.7649
$ws.Range("M134").Value = -2896.7649

$ws.Range("H136").Value = 14121.115
$ws.Range("I136").Value = 1369.85
$ws.Range("K136").Value = 4109.549999999999
$ws.Range("M136").Value = -1559.549999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1958.5217
$ws.Range("J129").Value = 2203.9
$ws.Range("L129").Value = 6611.700000000001
$ws.Range("N129").Value = -16611.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 39987
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents() | Out-Null

$ws.Range("H80").Value = 18749.5
$ws.Range("J80").Value = 27500
$ws.Range("L80").Value = 27500
$ws.Range("N80").Value = -29496

$ws.Range("H83").Value = 18749.5
$ws.Range("J83").Value = 27500
$ws.Range("L83").Value = 137500
$ws.Range("N83").Value = -147484

$ws.Range("H102").Value = 32458.45
$ws.Range("I102").Value = 48427.77
$ws.Range("J102").Value = 2801.1428
$ws.Range("K102").Value = 48427.77
$ws.Range("L102").Value = 2801.1428
$ws.Range("M102").Value = -46805.77
$ws.Range("N102").Value = -6045.1428

$ws.Range("H126").Value = 3229.1
$ws.Range("I126").Value = 2613.7144
$ws.Range("J126").Value = 4665
$ws.Range("K126").Value = 7841.1432
$ws.Range("L126").Value = 13995
$ws.Range("M126").Value = -5371.1432
$ws.Range("N126").Value = -18935

$ws.Range("H132").Value = 1943.375
$ws.Range("I132").Value = 1627.2069
$ws.Range("K132").Value = 4881.620699999999
$ws.Range("M132").Value = -2351.620699999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3373333
$ws.Range("J2").Value = 3373333
$ws.Range("L2").Value = 3373333
$ws.Range("N2").Value = -3373557

$ws.Range("H13").Value = 8500
$ws.Range("I13").Value = 5000
$ws.Range("J13").Value = 12000
$ws.Range("K13").Value = 5000
$ws.Range("L13").Value = 12000
$ws.Range("M13").Value = -4860
$ws.Range("N13").Value = -12280

$ws.Range("H18").Value = 20000
$ws.Range("J18").Value = 20000
$ws.Range("L18").Value = 20000
$ws.Range("N18").Value = -20344

$ws.Range("H24").Value = 19000
$ws.Range("J24").Value = 19000
$ws.Range("L24").Value = 19000
$ws.Range("N24").Value = -19686

$ws.Range("H55").Value = 914.5625
$ws.Range("I55").Value = 176
$ws.Range("J55").Value = 2145.5
$ws.Range("K55").Value = 176
$ws.Range("L55").Value = 2145.5
$ws.Range("M55").Value = -3
$ws.Range("N55").Value = -2491.5

$ws.Range("H132").Value = 2045.6111
$ws.Range("I132").Value = 916.25
$ws.Range("K132").Value = 2748.75
$ws.Range("M132").Value = -218.75

$ws.Range("H136").Value = 4196.3
$ws.Range("I136").Value = 3867.1304
$ws.Range("J136").Value = 5277.857
$ws.Range("K136").Value = 11601.3912
$ws.Range("L136").Value = 15833.571
$ws.Range("M136").Value = -9051.3912
$ws.Range("N136").Value = -20933.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2624.75
$ws.Range("J113").Value = 3349.5
$ws.Range("L113").Value = 10048.5
$ws.Range("N113").Value = -14388.5

$ws.Range("H122").Value = 463003
$ws.Range("I122").Value = 463003
$ws.Range("K122").Value = 1389009
$ws.Range("M122").Value = -1386559

$ws.Range("H126").Value = 175417.03
$ws.Range("I126").Value = 2713.6191
$ws.Range("K126").Value = 8140.8573
$ws.Range("M126").Value = -5670.8573

$ws.Range("H132").Value = 1615.9656
$ws.Range("I132").Value = 1198.8096
$ws.Range("K132").Value = 3596.4288
$ws.Range("M132").Value = -1066.4288
